# Revert "adding term 2.0 now utf-8"
# Restores the workbook to the prior 1.1.0 metadata and removes the
# "d7ff926a-4955-478f-b300-0b0ec0785013" concept row added on the
# "Include from FSIII" sheet.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet: revert Version / Date / Contact / Description ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value  = "1.1.0"
$meta.Range("B8").Value  = "2023-07-10T23:08:03+02:00"
$meta.Range("B10").Value = "No display for ContactDetail"
$meta.Range("B11").Value = "Matter of interes values to support when no observations have been made"

# --- "Include from FSIII" sheet: drop the row for the newly-added concept ---
$fsiii = $wb.Worksheets.Item("Include from FSIII")
$fsiii.Rows(2).Delete()
